# Actualización automática del inventario: se agrega un nuevo producto
# (Cartucho Damper Epson) como fila 43 de la hoja de inventario.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 43

$ws.Cells.Item($row, 1).Value = "XKXDN3"
$ws.Cells.Item($row, 2).Value = "Cartucho Damper Epson"
$ws.Cells.Item($row, 3).Value = "L1110 L1210 L1250 L3100 L3101 L3110 L3150 L3151 L3160 L3210 L3250 L4150 L4160 L4167 L4260 L5190 L5290 L5590 L6160 L6161 L6170 L6171 L6190 L6191 L6290 ET3750 ET4750"
$ws.Cells.Item($row, 4).Value = 25000
$ws.Cells.Item($row, 5).Value = 75000
$ws.Cells.Item($row, 6).Value = 20
$ws.Cells.Item($row, 7).Value = 30
$ws.Cells.Item($row, 8).Formula = "=(E43-D43)*G43"
$ws.Cells.Item($row, 9).Formula = "=D43*F43"
$ws.Cells.Item($row, 10).Value = 500000
